$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (36 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 23616
$ws.Range("I64").Value = 4293.3335
$ws.Range("J64").Value = 52600
$ws.Range("K64").Value = 4293.3335
$ws.Range("L64").Value = 52600
$ws.Range("M64").Value = -4045.3335
$ws.Range("N64").Value = -53096
$ws.Range("H67").Value = 23616
$ws.Range("I67").Value = 4293.3335
$ws.Range("J67").Value = 52600
$ws.Range("K67").Value = 4293.3335
$ws.Range("L67").Value = 52600
$ws.Range("M67").Value = -3435.3335
$ws.Range("N67").Value = -54316
$ws.Range("H125").Value = 5950
$ws.Range("I125").Value = 8900
$ws.Range("K125").Value = 80100
$ws.Range("M125").Value = -77640
$ws.Range("H137").Value = 1524.7778
$ws.Range("I137").Value = 1309.4375
$ws.Range("K137").Value = 3928.3125
$ws.Range("M137").Value = -1378.3125
$ws.Range("H138").Value = 3774.7544
$ws.Range("I138").Value = 3067.8928
$ws.Range("J138").Value = 4457.241
$ws.Range("K138").Value = 9203.678400000001
$ws.Range("L138").Value = 13371.723
$ws.Range("M138").Value = -4063.678400000001
$ws.Range("N138").Value = -23651.723
$ws.Range("H141").Value = 12200598
$ws.Range("I141").Value = 14709949
$ws.Range("J141").Value = 12319.714
$ws.Range("K141").Value = 44129847
$ws.Range("L141").Value = 36959.142
$ws.Range("M141").Value = -44124667
$ws.Range("N141").Value = -47319.142

# --- Sheet: ARM (16 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4310.2
$ws.Range("I32").Value = 4273.921
$ws.Range("K32").Value = 4273.921
$ws.Range("M32").Value = -3986.921
$ws.Range("H74").Value = 1267.3043
$ws.Range("I74").Value = 1032.45
$ws.Range("K74").Value = 1032.45
$ws.Range("M74").Value = -158.45
$ws.Range("H77").Value = 1267.3043
$ws.Range("I77").Value = 1032.45
$ws.Range("K77").Value = 5162.25
$ws.Range("M77").Value = -794.25
$ws.Range("H97").Value = 1694.8096
$ws.Range("J97").Value = 3199.5715
$ws.Range("L97").Value = 3199.5715
$ws.Range("N97").Value = -4191.5715

# --- Sheet: BSM (21 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1360.1482
$ws.Range("I20").Value = 1375.7894
$ws.Range("J20").Value = 1323
$ws.Range("K20").Value = 1375.7894
$ws.Range("L20").Value = 1323
$ws.Range("M20").Value = -1128.7894
$ws.Range("N20").Value = -1817
$ws.Range("H86").Value = 4790.75
$ws.Range("I86").Value = 1648.1666
$ws.Range("J86").Value = 7933.3335
$ws.Range("K86").Value = 1648.1666
$ws.Range("L86").Value = 7933.3335
$ws.Range("M86").Value = -525.1666
$ws.Range("N86").Value = -10179.3335
$ws.Range("H89").Value = 4790.75
$ws.Range("I89").Value = 1648.1666
$ws.Range("J89").Value = 7933.3335
$ws.Range("K89").Value = 8240.833000000001
$ws.Range("L89").Value = 39666.6675
$ws.Range("M89").Value = -2624.833000000001
$ws.Range("N89").Value = -50898.6675

# --- Sheet: CRP (18 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1920.5883
$ws.Range("I31").Value = 1677.8975
$ws.Range("J31").Value = 2709.3333
$ws.Range("K31").Value = 1677.8975
$ws.Range("L31").Value = 2709.3333
$ws.Range("M31").Value = -1382.8975
$ws.Range("N31").Value = -3299.3333
$ws.Range("H34").Value = 1920.5883
$ws.Range("I34").Value = 1677.8975
$ws.Range("J34").Value = 2709.3333
$ws.Range("K34").Value = 1677.8975
$ws.Range("L34").Value = 2709.3333
$ws.Range("M34").Value = -1475.8975
$ws.Range("N34").Value = -3113.3333
$ws.Range("H122").Value = 2429.0908
$ws.Range("I122").Value = 1934.6
$ws.Range("K122").Value = 5803.799999999999
$ws.Range("M122").Value = -3353.799999999999

# --- Sheet: CUL (32 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 711.0476
$ws.Range("I38").Value = 791.25
$ws.Range("J38").Value = 692.17645
$ws.Range("K38").Value = 2373.75
$ws.Range("L38").Value = 2076.52935
$ws.Range("M38").Value = -2026.75
$ws.Range("N38").Value = -2770.52935
$ws.Range("H68").Value = 1766.8334
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 1855.7778
$ws.Range("K68").Value = 4500
$ws.Range("L68").Value = 5567.3334
$ws.Range("M68").Value = -3689
$ws.Range("N68").Value = -7189.3334
$ws.Range("H71").Value = 1766.8334
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 1855.7778
$ws.Range("K71").Value = 13500
$ws.Range("L71").Value = 16702.0002
$ws.Range("M71").Value = -9444
$ws.Range("N71").Value = -24814.0002
$ws.Range("H107").Value = 2337064.5
$ws.Range("J107").Value = 3374494.2
$ws.Range("L107").Value = 10123482.6
$ws.Range("N107").Value = -10127322.6
$ws.Range("H113").Value = 1534.8572
$ws.Range("I113").Value = 497
$ws.Range("J113").Value = 1707.8334
$ws.Range("K113").Value = 1491
$ws.Range("L113").Value = 5123.5002
$ws.Range("M113").Value = 679
$ws.Range("N113").Value = -9463.5002

# --- Sheet: GSM (28 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6075.05
$ws.Range("I70").Value = 4443.4443
$ws.Range("J70").Value = 7410
$ws.Range("K70").Value = 4443.4443
$ws.Range("L70").Value = 7410
$ws.Range("M70").Value = -4173.4443
$ws.Range("N70").Value = -7950
$ws.Range("H73").Value = 6075.05
$ws.Range("I73").Value = 4443.4443
$ws.Range("J73").Value = 7410
$ws.Range("K73").Value = 4443.4443
$ws.Range("L73").Value = 7410
$ws.Range("M73").Value = -3507.4443
$ws.Range("N73").Value = -9282
$ws.Range("H97").Value = 725.17645
$ws.Range("I97").Value = 450.5
$ws.Range("J97").Value = 2007
$ws.Range("K97").Value = 450.5
$ws.Range("L97").Value = 2007
$ws.Range("M97").Value = 45.5
$ws.Range("N97").Value = -2999
$ws.Range("H113").Value = 1091745
$ws.Range("I113").Value = 2541.8333
$ws.Range("J113").Value = 1685855.8
$ws.Range("K113").Value = 2541.8333
$ws.Range("L113").Value = 1685855.8
$ws.Range("M113").Value = -371.8332999999998
$ws.Range("N113").Value = -1690195.8

# --- Sheet: LTW (14 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4169243.8
$ws.Range("I68").Value = 10418009
$ws.Range("J68").Value = 3400.3333
$ws.Range("K68").Value = 10418009
$ws.Range("L68").Value = 3400.3333
$ws.Range("M68").Value = -10417260
$ws.Range("N68").Value = -4898.3333
$ws.Range("H71").Value = 4169243.8
$ws.Range("I71").Value = 10418009
$ws.Range("J71").Value = 3400.3333
$ws.Range("K71").Value = 52090045
$ws.Range("L71").Value = 17001.6665
$ws.Range("M71").Value = -52086301
$ws.Range("N71").Value = -24489.6665

# --- Sheet: WVR (22 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7481.727
$ws.Range("I62").Value = 3800
$ws.Range("J62").Value = 9585.571
$ws.Range("K62").Value = 3800
$ws.Range("L62").Value = 9585.571
$ws.Range("M62").Value = -3176
$ws.Range("N62").Value = -10833.571
$ws.Range("H65").Value = 7481.727
$ws.Range("I65").Value = 3800
$ws.Range("J65").Value = 9585.571
$ws.Range("K65").Value = 19000
$ws.Range("L65").Value = 47927.855
$ws.Range("M65").Value = -15880
$ws.Range("N65").Value = -54167.855
$ws.Range("H112").Value = 27999.5
$ws.Range("J112").Value = 27999.5
$ws.Range("L112").Value = 27999.5
$ws.Range("N112").Value = -30953.5
$ws.Range("H113").Value = 1436.8125
$ws.Range("I113").Value = 1262.2727
$ws.Range("K113").Value = 3786.8181
$ws.Range("M113").Value = -1616.8181

Write-Host "Applied 187 cell updates across 8 sheets"
